$d = $word.ActiveDocument

$replacements = @(
    @{old = "346×8="; new = "641×8="},
    @{old = "225×8="; new = "657×9="},
    @{old = "874×3="; new = "772×2="},
    @{old = "612×9="; new = "997×2="},
    @{old = "640×9="; new = "353×2="},
    @{old = "564×6="; new = "129×8="},
    @{old = "201×5="; new = "114×3="},
    @{old = "366×2="; new = "507×3="},
    @{old = "548×8="; new = "538×4="},
    @{old = "866×5="; new = "291×4="},
    @{old = "527×2="; new = "652×7="},
    @{old = "435×9="; new = "619×5="},
    @{old = "261×3="; new = "429×5="},
    @{old = "607×8="; new = "580×5="},
    @{old = "876×4="; new = "950×5="},
    @{old = "112×7="; new = "393×9="},
    @{old = "476×7="; new = "932×5="},
    @{old = "333×5="; new = "930×8="},
    @{old = "838×7="; new = "463×3="},
    @{old = "132×2="; new = "196×2="},
    @{old = "868×6="; new = "136×5="},
    @{old = "430×3="; new = "534×9="},
    @{old = "590×8="; new = "643×8="},
    @{old = "610×9="; new = "457×6="},
    @{old = "908×3="; new = "326×2="}
)

foreach ($r in $replacements) {
    $range = $d.Content
    $range.Find.Execute($r.old, $true, $false, $false, $false, $false, $true, 1, $false, $r.new, 2)
}
